$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column-D price cells that look numeric to stay as text, matching the
# original inline-string content (e.g. "5.230", "1.0000", "37.00").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.905.80"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.862.69"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.85"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5052"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3624"
$ws.Range("E8").Value = "  -3.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07171"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8958"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.71"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07470"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.822.33"
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.89"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.230"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008469"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.16"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.940.07"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.028"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.083.03"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.34"
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.86"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.88"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.057"
$ws.Range("E28").Value = "  -2.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.08"
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.679"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09265"
$ws.Range("E32").Value = "  +2.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05084"
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.001"
$ws.Range("E34").Value = "  -3.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7435"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.280"
$ws.Range("E37").Value = "  +7.74%  "
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5612"
$ws.Range("E38").Value = "  +4.60%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02000"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.503"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "118.78"
$ws.Range("E42").Value = "  +2.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.474"
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.514"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1468"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4742"
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.0000"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.04"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.563"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.00"
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.00"
$ws.Range("E51").Value = "  -2.48%  "
